$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Update postStimBlankT (column J) from 250 to 500 for all data rows ---
for ($r = 2; $r -le 12; $r++) {
    $ws.Cells.Item($r, 10).Value2 = 500
}

# --- Add new columns O:S with headers ---
# Write in this order so the shared-string table grows as:
#   maskOnOff, nRevs, stairUp, stairDn, dnDivUp
$ws.Range("O1").Value2 = "maskOnOff"
$ws.Range("S1").Value2 = "nRevs"
$ws.Range("P1").Value2 = "stairUp"
$ws.Range("Q1").Value2 = "stairDn"
$ws.Range("R1").Value2 = "dnDivUp"

# --- Fill new column data for each data row ---
for ($r = 2; $r -le 12; $r++) {
    $ws.Cells.Item($r, 15).Value2 = 1   # O: maskOnOff
    $ws.Cells.Item($r, 16).Value2 = 1   # P: stairUp
    $ws.Cells.Item($r, 17).Value2 = 1   # Q: stairDn
    $ws.Cells.Item($r, 18).Value2 = 1   # R: dnDivUp
    $ws.Cells.Item($r, 19).Value2 = 12  # S: nRevs
}

# --- Update selection to reflect the new active cell ---
$ws.Range("J1").Select()
